$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 286152.5
$ws.Range("J17").Value = 303417.4
$ws.Range("L17").Value = 910252.2000000001
$ws.Range("N17").Value = -910588.2000000001

$ws.Range("H70").Value = 9476.471
$ws.Range("J70").Value = 11930.77
$ws.Range("L70").Value = 35792.31
$ws.Range("N70").Value = -36332.31

$ws.Range("H73").Value = 9476.471
$ws.Range("J73").Value = 11930.77
$ws.Range("L73").Value = 35792.31
$ws.Range("N73").Value = -37664.31

$ws.Range("H107").Value = 31253312
$ws.Range("I107").Value = 20836748
$ws.Range("K107").Value = 20836748
$ws.Range("M107").Value = -20834828

$ws.Range("H109").Value = 45000
$ws.Range("I109").Value = 45000
$ws.Range("K109").Value = 45000
$ws.Range("M109").Value = -43613

$ws.Range("H110").Value = 79999.5
$ws.Range("I110").Value = 80000
$ws.Range("K110").Value = 80000
$ws.Range("M110").Value = -75910

$ws.Range("H111").Value = 992.0909
$ws.Range("I111").Value = 951.3
$ws.Range("J111").Value = 1400
$ws.Range("K111").Value = 2853.9
$ws.Range("L111").Value = 4200
$ws.Range("M111").Value = 213.1000000000004
$ws.Range("N111").Value = -10334

$ws.Range("H113").Value = 3900
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = 554

$ws.Range("H115").Value = 6508220
$ws.Range("I115").Value = 7808864.5
$ws.Range("K115").Value = 23426593.5
$ws.Range("M115").Value = -23425026.5

$ws.Range("H116").Value = 6666
$ws.Range("I116").Value = 4999
$ws.Range("K116").Value = 4999
$ws.Range("M116").Value = -1557

$ws.Range("H118").Value = 4080369.2
$ws.Range("I118").Value = 4080369.2
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 12241107.6
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -12239450.6

$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676

$ws.Range("H131").Value = 3798.1765
$ws.Range("I131").Value = 1904.6666
$ws.Range("K131").Value = 5713.9998
$ws.Range("M131").Value = -673.9997999999996

$ws.Range("H132").Value = 6711.2104
$ws.Range("I132").Value = 6146.2144
$ws.Range("K132").Value = 18438.6432
$ws.Range("M132").Value = -15908.6432

$ws.Range("H137").Value = 31472.7
$ws.Range("I137").Value = 48166.625
$ws.Range("J137").Value = 6431.8125
$ws.Range("K137").Value = 144499.875
$ws.Range("L137").Value = 19295.4375
$ws.Range("M137").Value = -141949.875
$ws.Range("N137").Value = -24395.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1697
$ws.Range("I45").Value = 1748
$ws.Range("K45").Value = 1748
$ws.Range("M45").Value = -1371

$ws.Range("H88").Value = 1940.625
$ws.Range("I88").Value = 1100
$ws.Range("J88").Value = 2220.8333
$ws.Range("K88").Value = 1100
$ws.Range("L88").Value = 2220.8333
$ws.Range("M88").Value = -694
$ws.Range("N88").Value = -3032.8333

$ws.Range("H91").Value = 1940.625
$ws.Range("I91").Value = 1100
$ws.Range("J91").Value = 2220.8333
$ws.Range("K91").Value = 1100
$ws.Range("L91").Value = 2220.8333
$ws.Range("M91").Value = 304
$ws.Range("N91").Value = -5028.8333

$ws.Range("H122").Value = 2564.4211
$ws.Range("I122").Value = 1787.3846
$ws.Range("K122").Value = 5362.1538
$ws.Range("M122").Value = -2912.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 1626.4
$ws.Range("I10").Value = 2675
$ws.Range("K10").Value = 2675
$ws.Range("M10").Value = -2535

$ws.Range("H20").Value = 1097.76
$ws.Range("I20").Value = 973.2
$ws.Range("K20").Value = 973.2
$ws.Range("M20").Value = -726.2

$ws.Range("H86").Value = 5895.885
$ws.Range("I86").Value = 5818.8096
$ws.Range("K86").Value = 5818.8096
$ws.Range("M86").Value = -4695.8096

$ws.Range("H89").Value = 5895.885
$ws.Range("I89").Value = 5818.8096
$ws.Range("K89").Value = 29094.048
$ws.Range("M89").Value = -23478.048

$ws.Range("H94").Value = 1015
$ws.Range("I94").Value = 1058.3125
$ws.Range("J94").Value = 899.5
$ws.Range("K94").Value = 1058.3125
$ws.Range("L94").Value = 899.5
$ws.Range("M94").Value = -607.3125
$ws.Range("N94").Value = -1801.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 4980
$ws.Range("J11").Value = 4980
$ws.Range("L11").Value = 4980
$ws.Range("N11").Value = -5260

$ws.Range("H134").Value = 7880.522
$ws.Range("I134").Value = 9250.764999999999
$ws.Range("J134").Value = 3998.1667
$ws.Range("K134").Value = 27752.295
$ws.Range("L134").Value = 11994.5001
$ws.Range("M134").Value = -25217.295
$ws.Range("N134").Value = -17064.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2112.25
$ws.Range("I41").Value = 99
$ws.Range("J41").Value = 2783.3333
$ws.Range("K41").Value = 297
$ws.Range("L41").Value = 8349.999899999999
$ws.Range("M41").Value = 41
$ws.Range("N41").Value = -9025.999899999999

$ws.Range("H98").Value = 355.81818
$ws.Range("I98").Value = 333
$ws.Range("J98").Value = 368.85715
$ws.Range("K98").Value = 999
$ws.Range("L98").Value = 1106.57145
$ws.Range("M98").Value = 499
$ws.Range("N98").Value = -4102.571449999999

$ws.Range("H107").Value = 527.05884
$ws.Range("J107").Value = 584
$ws.Range("L107").Value = 1752
$ws.Range("N107").Value = -5592

$ws.Range("H108").Value = 38
$ws.Range("I108").Value = 38
$ws.Range("K108").Value = 114
$ws.Range("M108").Value = 2766

$ws.Range("H109").Value = 1000000
$ws.Range("I109").Value = 1000000
$ws.Range("K109").Value = 3000000
$ws.Range("M109").Value = -2998960

$ws.Range("H111").Value = 499.5
$ws.Range("I111").Value = 499.5
$ws.Range("K111").Value = 1498.5
$ws.Range("M111").Value = 1568.5

$ws.Range("H113").Value = 593.7727
$ws.Range("I113").Value = 342.85715
$ws.Range("J113").Value = 710.86664
$ws.Range("K113").Value = 1028.57145
$ws.Range("L113").Value = 2132.59992
$ws.Range("M113").Value = 1141.42855
$ws.Range("N113").Value = -6472.59992

$ws.Range("H114").Value = 637.8570999999999
$ws.Range("I114").Value = 466.8
$ws.Range("J114").Value = 1065.5
$ws.Range("K114").Value = 1400.4
$ws.Range("L114").Value = 3196.5
$ws.Range("M114").Value = 1853.6
$ws.Range("N114").Value = -9704.5

$ws.Range("H116").Value = 21077.6
$ws.Range("I116").Value = 26064
$ws.Range("J116").Value = 1132
$ws.Range("K116").Value = 78192
$ws.Range("L116").Value = 3396
$ws.Range("M116").Value = -74750
$ws.Range("N116").Value = -10280

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H119").Value = 1500
$ws.Range("I119").Value = 1500
$ws.Range("K119").Value = 4500
$ws.Range("M119").Value = 338

$ws.Range("H121").Value = 148599
$ws.Range("I121").Value = 4998.5
$ws.Range("K121").Value = 14995.5
$ws.Range("M121").Value = -13685.5

$ws.Range("H122").Value = 168
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 1900
$ws.Range("J123").Value = 1900
$ws.Range("L123").Value = 5700
$ws.Range("N123").Value = -10600

$ws.Range("H125").Value = 1450
$ws.Range("I125").Value = 1450
$ws.Range("K125").Value = 4350
$ws.Range("M125").Value = 570

$ws.Range("H129").Value = 25642082
$ws.Range("I129").Value = 41667172
$ws.Range("J129").Value = 1933.2
$ws.Range("K129").Value = 125001516
$ws.Range("L129").Value = 5799.6
$ws.Range("M129").Value = -124996516
$ws.Range("N129").Value = -15799.6

$ws.Range("H131").Value = 8334759
$ws.Range("I131").Value = 83334250
$ws.Range("J131").Value = 1482.1852
$ws.Range("K131").Value = 250002750
$ws.Range("L131").Value = 4446.5556
$ws.Range("M131").Value = -249997710
$ws.Range("N131").Value = -14526.5556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 25499.666
$ws.Range("I99").Value = 25499.666
$ws.Range("K99").Value = 25499.666
$ws.Range("M99").Value = -23253.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1198.1538
$ws.Range("I16").Value = 679.9091
$ws.Range("K16").Value = 679.9091
$ws.Range("M16").Value = -509.9091

$ws.Range("H46").Value = 3275.4666
$ws.Range("I46").Value = 2814.1
$ws.Range("K46").Value = 2814.1
$ws.Range("M46").Value = -2626.1

$ws.Range("H107").Value = 172322.33
$ws.Range("I107").Value = 172322.33
$ws.Range("K107").Value = 172322.33
$ws.Range("M107").Value = -170402.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 821.8333
$ws.Range("I100").Value = 1130.2858
$ws.Range("J100").Value = 390
$ws.Range("K100").Value = 2260.5716
$ws.Range("L100").Value = 780
$ws.Range("M100").Value = -1719.5716
$ws.Range("N100").Value = -1862

$ws.Range("H107").Value = 612.3333
$ws.Range("I107").Value = 667.6
$ws.Range("K107").Value = 2002.8
$ws.Range("M107").Value = -82.80000000000018

$ws.Range("H138").Value = 74332.336
$ws.Range("J138").Value = 74332.336
$ws.Range("L138").Value = 74332.336
$ws.Range("N138").Value = -84612.336

$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360
